$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B1 is stored as literal text (inline string), not a date serial number.
$ws.Range("B1").Value = "14/03/2023"

$ws.Range("B2").Value = 1008.8
$ws.Range("C2").Value = 76

$ws.Range("B3").Value = 120
$ws.Range("C3").Value = 76

$ws.Range("B4").Value = 112

$ws.Range("B5").Value = 678
$ws.Range("C5").Value = 76

$ws.Range("B6").Value = 536
$ws.Range("C6").Value = 76

$ws.Range("B7").Value = 165
$ws.Range("C7").Value = 76

$ws.Range("B8").Value = 124

$ws.Range("B9").Value = 591
$ws.Range("C9").Value = 76

$ws.Range("B10").Value = 86
$ws.Range("C10").Value = 38

$ws.Range("C11").Value = 38

$ws.Range("B12").Value = 68
$ws.Range("C12").Value = 76
